$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 7 (IGLUW member) - JUEVES / VIERNES answers
$ws.Range("F7").Value = "Se cumplio lo planeado"
$ws.Range("G7").Value = "Se cumplio lo propuesto"

# Row 8 - JUEVES / VIERNES answers
$ws.Range("F8").Value = "Realizar la presentación sobre los avances obtenidos y toma de sugerencias"
$ws.Range("G8").Value = "Asistir a la reunión para la asignación de tareas y tener en cuenta qué dijo el profe."

# Row 9 - JUEVES / VIERNES answers
$ws.Range("F9").Value = "Nada"
$ws.Range("G9").Value = "Nada"

# Update selection / view to match
[void]$ws.Range("G8").Select()
